$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at G, shifting Patient#, Cell#, Platform, DataID one column right.
$ws.Range("G1").EntireColumn.Insert()
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# New "Procedure" header and values
$ws.Range("G1").Value = "Procedure"
$ws.Range("G2").Value = "Transplant"
$ws.Range("G3").Value = "Endarterectomy"
$ws.Range("G4").Value = "Endarterectomy"
$ws.Range("G3:G4").Style = "Normal"

# Fix author name typo
$ws.Range("A3").Value = "Slenders et al."

# Update tissue descriptions
$ws.Range("F2").Value = "Coronary Artery"
$ws.Range("F3").Value = "Carotid Artery Plaque"
$ws.Range("F4").Value = "Carotid Artery Plaque"

# Update selection to match final state
$ws.Range("A4").Select()
